# Scheduled runner update: refresh Leve profit-margin figures (currentAveragePrice /
# NQ / HQ prices and computed leve profits) across all job sheets following the
# latest Universalis market-board pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1366.6666
$ws.Range("J70").Value = 1100
$ws.Range("L70").Value = 3300
$ws.Range("N70").Value = -3840
$ws.Range("H73").Value = 1366.6666
$ws.Range("J73").Value = 1100
$ws.Range("L73").Value = 3300
$ws.Range("N73").Value = -5172
$ws.Range("H116").Value = 17645294
$ws.Range("J116").Value = 5436.5454
$ws.Range("L116").Value = 5436.5454
$ws.Range("N116").Value = -12320.5454
$ws.Range("H129").Value = 295880.12
$ws.Range("J129").Value = 324487.97
$ws.Range("L129").Value = 973463.9099999999
$ws.Range("N129").Value = -983463.9099999999
$ws.Range("H132").Value = 34486090
$ws.Range("I132").Value = 38465120
$ws.Range("K132").Value = 115395360
$ws.Range("M132").Value = -115392830
$ws.Range("H135").Value = 13892249
$ws.Range("I135").Value = 625.25
$ws.Range("J135").Value = 125025240
$ws.Range("K135").Value = 5627.25
$ws.Range("L135").Value = 1125227160
$ws.Range("M135").Value = -3092.25
$ws.Range("N135").Value = -1125232230
$ws.Range("H137").Value = 67885.766
$ws.Range("I137").Value = 96337.88
$ws.Range("J137").Value = 1497.5
$ws.Range("K137").Value = 289013.64
$ws.Range("L137").Value = 4492.5
$ws.Range("M137").Value = -286463.64
$ws.Range("N137").Value = -9592.5
$ws.Range("H138").Value = 3459.5
$ws.Range("I138").Value = 2195.389
$ws.Range("J138").Value = 4334.654
$ws.Range("K138").Value = 6586.167
$ws.Range("L138").Value = 13003.962
$ws.Range("M138").Value = -1446.167
$ws.Range("N138").Value = -23283.962

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8774.562
$ws.Range("I32").Value = 6350.75
$ws.Range("K32").Value = 6350.75
$ws.Range("M32").Value = -6063.75
$ws.Range("H74").Value = 28572914
$ws.Range("I74").Value = 40000628
$ws.Range("K74").Value = 40000628
$ws.Range("M74").Value = -39999754
$ws.Range("H77").Value = 28572914
$ws.Range("I77").Value = 40000628
$ws.Range("K77").Value = 200003140
$ws.Range("M77").Value = -199998772
$ws.Range("H110").Value = 1079.7858
$ws.Range("J110").Value = 800
$ws.Range("L110").Value = 800
$ws.Range("N110").Value = -4890
$ws.Range("H122").Value = 2300.7827
$ws.Range("I122").Value = 2176.95
$ws.Range("K122").Value = 6530.849999999999
$ws.Range("M122").Value = -4080.849999999999
$ws.Range("H132").Value = 10216503
$ws.Range("I132").Value = 12501951
$ws.Range("K132").Value = 37505853
$ws.Range("M132").Value = -37503323

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1353443.6
$ws.Range("I105").Value = 1447.2307
$ws.Range("J105").Value = 2085775
$ws.Range("K105").Value = 1447.2307
$ws.Range("L105").Value = 2085775
$ws.Range("M105").Value = 299.7692999999999
$ws.Range("N105").Value = -2089269

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 904
$ws.Range("J16").Value = 680
$ws.Range("L16").Value = 680
$ws.Range("N16").Value = -1254
$ws.Range("H22").Value = 181.47058
$ws.Range("J22").Value = 240
$ws.Range("L22").Value = 240
$ws.Range("N22").Value = -940
$ws.Range("H31").Value = 4189.1113
$ws.Range("I31").Value = 2273.2
$ws.Range("J31").Value = 5840.759
$ws.Range("K31").Value = 2273.2
$ws.Range("L31").Value = 5840.759
$ws.Range("M31").Value = -1978.2
$ws.Range("N31").Value = -6430.759
$ws.Range("H34").Value = 4189.1113
$ws.Range("I34").Value = 2273.2
$ws.Range("J34").Value = 5840.759
$ws.Range("K34").Value = 2273.2
$ws.Range("L34").Value = 5840.759
$ws.Range("M34").Value = -2071.2
$ws.Range("N34").Value = -6244.759
$ws.Range("H62").Value = 2706.3777
$ws.Range("I62").Value = 2527.05
$ws.Range("J62").Value = 4141
$ws.Range("K62").Value = 2527.05
$ws.Range("L62").Value = 4141
$ws.Range("M62").Value = -1903.05
$ws.Range("N62").Value = -5389
$ws.Range("H65").Value = 2706.3777
$ws.Range("I65").Value = 2527.05
$ws.Range("J65").Value = 4141
$ws.Range("K65").Value = 12635.25
$ws.Range("L65").Value = 20705
$ws.Range("M65").Value = -9515.25
$ws.Range("N65").Value = -26945
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H97").Value = 31999
$ws.Range("J97").Value = 31999
$ws.Range("L97").Value = 31999
$ws.Range("N97").Value = -33981
$ws.Range("H99").Value = 3627.7693
$ws.Range("I99").Value = 2698.6
$ws.Range("K99").Value = 2698.6
$ws.Range("M99").Value = -1200.6
$ws.Range("H105").Value = 2845.7144
$ws.Range("I105").Value = 2781.8
$ws.Range("K105").Value = 2781.8
$ws.Range("M105").Value = -1034.8
$ws.Range("H113").Value = 904
$ws.Range("J113").Value = 680
$ws.Range("L113").Value = 680
$ws.Range("N113").Value = -5020
$ws.Range("H126").Value = 3627.7693
$ws.Range("I126").Value = 2698.6
$ws.Range("K126").Value = 8095.799999999999
$ws.Range("M126").Value = -5625.799999999999
$ws.Range("H132").Value = 41669316
$ws.Range("I132").Value = 47620780
$ws.Range("K132").Value = 142862340
$ws.Range("M132").Value = -142859810

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 367
$ws.Range("I23").Value = 20
$ws.Range("J23").Value = 398.54544
$ws.Range("K23").Value = 60
$ws.Range("L23").Value = 1195.63632
$ws.Range("M23").Value = 175
$ws.Range("N23").Value = -1665.63632
$ws.Range("H131").Value = 689.67
$ws.Range("J131").Value = 735.83905
$ws.Range("L131").Value = 2207.51715
$ws.Range("N131").Value = -12287.51715
$ws.Range("H132").Value = 996
$ws.Range("J132").Value = 993.3333
$ws.Range("L132").Value = 8939.9997
$ws.Range("N132").Value = -13999.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3476583
$ws.Range("I70").Value = 4444.1816
$ws.Range("J70").Value = 8932801
$ws.Range("K70").Value = 4444.1816
$ws.Range("L70").Value = 8932801
$ws.Range("M70").Value = -4174.1816
$ws.Range("N70").Value = -8933341
$ws.Range("H73").Value = 3476583
$ws.Range("I73").Value = 4444.1816
$ws.Range("J73").Value = 8932801
$ws.Range("K73").Value = 4444.1816
$ws.Range("L73").Value = 8932801
$ws.Range("M73").Value = -3508.1816
$ws.Range("N73").Value = -8934673
$ws.Range("H126").Value = 3524
$ws.Range("I126").Value = 2465.2942
$ws.Range("K126").Value = 7395.882599999999
$ws.Range("M126").Value = -4925.882599999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4882.533
$ws.Range("I40").Value = 4476.1816
$ws.Range("K40").Value = 4476.1816
$ws.Range("M40").Value = -4340.1816
$ws.Range("H122").Value = 938715.6
$ws.Range("I122").Value = 1405302
$ws.Range("J122").Value = 5542.857
$ws.Range("K122").Value = 4215906
$ws.Range("L122").Value = 16628.571
$ws.Range("M122").Value = -4213456
$ws.Range("N122").Value = -21528.571
$ws.Range("H136").Value = 2382.724
$ws.Range("I136").Value = 2382.724
$ws.Range("K136").Value = 7148.172
$ws.Range("M136").Value = -4598.172

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5042.857
$ws.Range("I62").Value = 4750
$ws.Range("J62").Value = 5160
$ws.Range("K62").Value = 4750
$ws.Range("L62").Value = 5160
$ws.Range("M62").Value = -4126
$ws.Range("N62").Value = -6408
$ws.Range("H65").Value = 5042.857
$ws.Range("I65").Value = 4750
$ws.Range("J65").Value = 5160
$ws.Range("K65").Value = 23750
$ws.Range("L65").Value = 25800
$ws.Range("M65").Value = -20630
$ws.Range("N65").Value = -32040
$ws.Range("H75").Value = 24000
$ws.Range("J75").Value = 24000
$ws.Range("L75").Value = 24000
$ws.Range("N75").Value = -25872
$ws.Range("H78").Value = 24000
$ws.Range("J78").Value = 24000
$ws.Range("L78").Value = 72000
$ws.Range("N78").Value = -81360
$ws.Range("H86").Value = 30000
$ws.Range("J86").Value = 30000
$ws.Range("L86").Value = 30000
$ws.Range("N86").Value = -32246
$ws.Range("H89").Value = 30000
$ws.Range("J89").Value = 30000
$ws.Range("L89").Value = 150000
$ws.Range("N89").Value = -161232
$ws.Range("H122").Value = 1887.5
$ws.Range("I122").Value = 1887.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5662.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3212.5
$ws.Range("N122").ClearContents()
